$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = 25.02000000000047
$ws.Range("H2").Value = [double]"2.135044278125301e-16"
$ws.Range("I2").ClearContents()
$ws.Range("K2").Value = 48.54583869701901
$ws.Range("L2").Value = "[39.512866284888744, 57.57881110914927]"
$ws.Range("O2").Value = 1.79250031285904
$ws.Range("P2").Value = "[1.5912371198362703, 1.9937635058818106]"
$ws.Range("S2").Value = 60.92972919956009
$ws.Range("T2").Value = "[55.64771353770939, 66.21174486141078]"
$ws.Range("W2").Value = 17.8821621621625
$ws.Range("X2").Value = 17.08072072072104
$ws.Range("Y2").Value = 18.68360360360396

# --- Row 3 updates ---
$ws.Range("E3").Value = 22.81000000000013
$ws.Range("H3").Value = [double]"2.135044278125301e-16"
$ws.Range("K3").Value = 46.03541875182044
$ws.Range("L3").Value = "[35.83902386822543, 56.231813635415456]"
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = -2.792526803190928
$ws.Range("P3").Value = "[-3.018947895341544, -2.5661057110403114]"
$ws.Range("S3").Value = 62.61087463266746
$ws.Range("T3").Value = "[56.96734552917039, 68.25440373616453]"
$ws.Range("W3").Value = 10.13777777777784
$ws.Range("X3").Value = 9.315795795795848
$ws.Range("Y3").Value = 10.95975975975982

Write-Output "Edit applied"
